$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '46.188.70'
$ws.Range("E2").Value = '  +3.57%  '

$ws.Range("D3").Value = '2.452.40'
$ws.Range("E3").Value = '  +1.16%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '321.15'
$ws.Range("E5").Value = '  +2.67%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '104.87'
$ws.Range("E6").Value = '  +3.22%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.519'
$ws.Range("E7").Value = '  +1.24%  '

$ws.Range("E8").Value = '  -0.03%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.537'
$ws.Range("E9").Value = '  +4.13%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.95'
$ws.Range("E10").Value = '  +1.98%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0807'
$ws.Range("E11").Value = '  +0.69%  '

$ws.Range("E12").Value = '  -3.16%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.28'
$ws.Range("E13").Value = '  -3.48%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.08'
$ws.Range("E14").Value = '  +2.03%  '

$ws.Range("D15").Value = '2.833.54'
$ws.Range("E15").Value = '  +1.07%  '

$ws.Range("D16").Value = '2.463.42'
$ws.Range("E16").Value = '  +1.34%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.845'
$ws.Range("E17").Value = '  +1.25%  '

$ws.Range("D18").Value = '46.031.59'
$ws.Range("E18").Value = '  +3.57%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.59'
$ws.Range("E19").Value = '  +0.52%  '

$ws.Range("E20").Value = '  +0.45%  '

$ws.Range("D21").Value = '0.0₃0933'
$ws.Range("E21").Value = '  +0.46%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '71.22'
$ws.Range("E22").Value = '  +3.44%  '

$ws.Range("E23").Value = '  +4.78%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '247.42'
$ws.Range("E24").Value = '  +2.59%  '

$ws.Range("E25").Value = '  +1.79%  '

$ws.Range("B26").Value = 'EthereumClassic'
$ws.Range("C26").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '25.93'
$ws.Range("E26").Value = '  +2.94%  '

$ws.Range("B27").Value = 'Dai'
$ws.Range("C27").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("E27").Value = '  +0.02%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.26'
$ws.Range("E28").Value = '  -0.20%  '

$ws.Range("E29").Value = '  +0.79%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '33.78'
$ws.Range("E30").Value = '  +1.16%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '49.37'
$ws.Range("E31").Value = '  +1.52%  '

$ws.Range("E32").Value = '  +2.92%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '20.01'
$ws.Range("E33").Value = '  +2.32%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.35'
$ws.Range("E34").Value = '  +3.15%  '

$ws.Range("E35").Value = '  -0.05%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0761'
$ws.Range("E36").Value = '  -1.04%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.54'
$ws.Range("E37").Value = '  +0.06%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.90'
$ws.Range("E38").Value = '  +0.39%  '

$ws.Range("E39").Value = '  +1.74%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '126.27'
$ws.Range("E40").Value = '  +1.80%  '

$ws.Range("E41").Value = '  +2.14%  '

$ws.Range("E42").Value = '  +1.57%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '20.91'
$ws.Range("E43").Value = '  -1.92%  '

$ws.Range("E44").Value = '  +1.15%  '

$ws.Range("D45").Value = '1.968.74'
$ws.Range("E45").Value = '  +0.98%  '

$ws.Range("E46").Value = '  +1.29%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.07'
$ws.Range("E47").Value = '  -5.17%  '

$ws.Range("E48").Value = '  +12.46%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.13'

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.02'
$ws.Range("E50").Value = '  +7.85%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '77.86'
$ws.Range("E51").Value = '  +5.30%  '
